$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Note: this engine's ColumnWidth setter bakes in a fixed +5/6 char padding
# offset when serialising to the OOXML <col width> attribute, so we back it
# out here to land on the exact target widths (53 / 89 / 29).
$ws.Columns.Item(3).ColumnWidth = 53 - 5/6
$ws.Columns.Item(4).ColumnWidth = 89 - 5/6
$ws.Columns.Item(8).ColumnWidth = 29 - 5/6

# --- Data rows 2-11 (A..H) ---
$data = @(
  @("1327027","https://aiesec.org/opportunity/global-talent/1327027","EMEA Marketing and Communication Trainee - only eu","Zaventem, Belgio","No","1 applicant","3 - 6 Months","TerumoBCT"),
  @("1327026","https://aiesec.org/opportunity/global-talent/1327026","EMEA Marketing Intern - ONLY EU","Zaventem, Belgio","No","1 applicant","3 - 6 Months","TerumoBCT"),
  @("1326988","https://aiesec.org/opportunity/global-talent/1326988","English Instructor and marketing specialist","Zagazig, El-Hariry, Zagazig 1, Al-Sharqia Governorate, Egypt","No","1 applicant","9 - 12 Weeks","we care academy zag"),
  @("1326555","https://aiesec.org/opportunity/global-talent/1326555","Graphic desgin","Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt","No","0 applicants","9 - 12 Weeks","Egypt holiday travel"),
  @("1326434","https://aiesec.org/opportunity/global-talent/1326434","Pharmacovigilance Quality Trainee (ONLY EU)","Bruxelles, Belgio","No","22 applicants","6 - 18 Months","UCB"),
  @("1325059","https://aiesec.org/opportunity/global-talent/1325059","Interior designer","Cairo, Cairo Governorate, Egypt","No","14 applicants","6 - 18 Months","Kaian Arabi"),
  @("1324910","https://aiesec.org/opportunity/global-talent/1324910","Graphic designer","Cairo, Cairo Governorate, Egypt","No","8 applicants","3 - 6 Months","Transition Agency"),
  @("1324592","https://aiesec.org/opportunity/global-talent/1324592","Digital marketing","New Damietta City, Damietta El-Gadeeda City, New Damietta, Damietta Governorate, Egypt","No","13 applicants","3 - 6 Months","Business Haven Consultancy"),
  @("1321333","https://aiesec.org/opportunity/global-talent/1321333","Key Account Manager","Bogotá, Colombia","No","21 applicants","6 - 18 Months","SHIPTAINER COLOMBIA S.A.S."),
  @("1317035","https://aiesec.org/opportunity/global-talent/1317035","Business Development Champion","Cairo, Cairo Governorate, Egypt","No","39 applicants","9 - 12 Weeks","ACT Management Consulting")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = "'" + $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $r = $r + 1
}

# Row 6 (E6) used to be highlighted "Yes" (special style); now plain "No" with default style.
$ws.Range("E6").Style = "Normal"
